# Update loading_percent values for rows 2-25 (A=0..23) per new Case_3_141 (380 kV) run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.2010195329214
$ws.Range("C2").Value = 6.348597204040211
$ws.Range("D2").Value = 6.166271433259684
$ws.Range("F2").Value = 32.11997103907203
$ws.Range("G2").Value = 3.67069451656726
$ws.Range("I2").Value = 25.43644112348629
$ws.Range("K2").Value = 10.8551038974116
$ws.Range("M2").Value = 19.97926429944394
$ws.Range("N2").Value = 20.63038108250789

$ws.Range("B3").Value = 9.993223557450916
$ws.Range("C3").Value = 6.161188816825276
$ws.Range("D3").Value = 6.166756878286209
$ws.Range("F3").Value = 31.99526619901115
$ws.Range("G3").Value = 3.673347356561079
$ws.Range("I3").Value = 25.43717054731999
$ws.Range("K3").Value = 10.7191196820789
$ws.Range("M3").Value = 19.39140724492064
$ws.Range("N3").Value = 20.6781232241942

$ws.Range("B4").Value = 9.866404362153347
$ws.Range("C4").Value = 6.04521137976755
$ws.Range("D4").Value = 6.166941496036576
$ws.Range("F4").Value = 31.92621034295823
$ws.Range("G4").Value = 3.67506028193894
$ws.Range("I4").Value = 25.442385677547
$ws.Range("K4").Value = 10.6377101444263
$ws.Range("M4").Value = 19.02955670327156
$ws.Range("N4").Value = 20.70932067105695

$ws.Range("B5").Value = 9.814999142454043
$ws.Range("C5").Value = 5.997804687530105
$ws.Range("D5").Value = 6.166988485123845
$ws.Range("F5").Value = 31.89997529364077
$ws.Range("G5").Value = 3.675779531503692
$ws.Range("I5").Value = 25.4457082773189
$ws.Range("K5").Value = 10.6051006883781
$ws.Range("M5").Value = 18.88211829205501
$ws.Range("N5").Value = 20.72250760243267

$ws.Range("B6").Value = 9.806482332217591
$ws.Range("C6").Value = 5.98992661230943
$ws.Range("D6").Value = 6.166994589996095
$ws.Range("F6").Value = 31.89573451996176
$ws.Range("G6").Value = 3.675900246191846
$ws.Range("I6").Value = 25.44633226154515
$ws.Range("K6").Value = 10.5997212558923
$ws.Range("M6").Value = 18.85764478051637
$ws.Range("N6").Value = 20.72472589632973

$ws.Range("B7").Value = 9.865709873006564
$ws.Range("C7").Value = 6.0445725099235
$ws.Range("D7").Value = 6.166942243734405
$ws.Range("F7").Value = 31.92584879154419
$ws.Range("G7").Value = 3.675069895971947
$ws.Range("I7").Value = 25.44242564126
$ws.Range("K7").Value = 10.63726801818557
$ws.Range("M7").Value = 19.02756788882086
$ws.Range("N7").Value = 20.70949659642248

$ws.Range("B8").Value = 10.12926045935751
$ws.Range("C8").Value = 6.284215801857049
$ws.Range("D8").Value = 6.166462531713703
$ws.Range("F8").Value = 32.07542401225545
$ws.Range("G8").Value = 3.671591814113186
$ws.Range("I8").Value = 25.43570250263079
$ws.Range("K8").Value = 10.80780878288352
$ws.Range("M8").Value = 19.77690275452946
$ws.Range("N8").Value = 20.64645168098699

$ws.Range("B9").Value = 10.64867931171664
$ws.Range("C9").Value = 6.743535970582528
$ws.Range("D9").Value = 6.164608701653846
$ws.Range("F9").Value = 32.4275230295202
$ws.Range("G9").Value = 3.66543479417148
$ws.Range("I9").Value = 25.46038765782751
$ws.Range("K9").Value = 11.15692790845445
$ws.Range("M9").Value = 21.22900276320528
$ws.Range("N9").Value = 20.53776248996507

$ws.Range("B10").Value = 11.02735497457906
$ws.Range("C10").Value = 7.070328222081887
$ws.Range("D10").Value = 6.16267179046686
$ws.Range("F10").Value = 32.72077637995417
$ws.Range("G10").Value = 3.661310695073342
$ws.Range("I10").Value = 25.50164213003028
$ws.Range("K10").Value = 11.41985545504393
$ws.Range("M10").Value = 22.27269591385013
$ws.Range("N10").Value = 20.46701215344883

$ws.Range("B11").Value = 11.19808885430325
$ws.Range("C11").Value = 7.21590315108722
$ws.Range("D11").Value = 6.161661825856123
$ws.Range("F11").Value = 32.86137361679501
$ws.Range("G11").Value = 3.659520201811643
$ws.Range("I11").Value = 25.52542545473425
$ws.Range("K11").Value = 11.54032632733863
$ws.Range("M11").Value = 22.74012503800054
$ws.Range("N11").Value = 20.43680206305341

$ws.Range("B12").Value = 11.26245202444076
$ws.Range("C12").Value = 7.270528273706985
$ws.Range("D12").Value = 6.161260527194526
$ws.Range("F12").Value = 32.9156187062561
$ws.Range("G12").Value = 3.658854412477075
$ws.Range("I12").Value = 25.53515146059778
$ws.Range("K12").Value = 11.58602534284838
$ws.Range("M12").Value = 22.91589829343211
$ws.Range("N12").Value = 20.42564628467274

$ws.Range("B13").Value = 11.24860424031009
$ws.Range("C13").Value = 7.258786910389146
$ws.Range("D13").Value = 6.161347797354096
$ws.Range("F13").Value = 32.9038919116652
$ws.Range("G13").Value = 3.658997259306576
$ws.Range("I13").Value = 25.53302481037755
$ws.Range("K13").Value = 11.57618042071138
$ws.Range("M13").Value = 22.87809977694812
$ws.Range("N13").Value = 20.42803624132246

$ws.Range("B14").Value = 11.20339025232504
$ws.Range("C14").Value = 7.220407541368791
$ws.Range("D14").Value = 6.161629190331059
$ws.Range("F14").Value = 32.86581645412421
$ws.Range("G14").Value = 3.659465182210307
$ws.Range("I14").Value = 25.52621120485536
$ws.Range("K14").Value = 11.54408467753671
$ws.Range("M14").Value = 22.75461164242799
$ws.Range("N14").Value = 20.43587857499422

$ws.Range("B15").Value = 11.17565558880799
$ws.Range("C15").Value = 7.196832224012518
$ws.Range("D15").Value = 6.161799087721491
$ws.Range("F15").Value = 32.842623941658
$ws.Range("G15").Value = 3.659753389201812
$ws.Range("I15").Value = 25.52213135388769
$ws.Range("K15").Value = 11.52443411826782
$ws.Range("M15").Value = 22.67880624908973
$ws.Range("N15").Value = 20.44071923837658

$ws.Range("B16").Value = 11.01616026480084
$ws.Range("C16").Value = 7.060747455504951
$ws.Range("D16").Value = 6.162735174052234
$ws.Range("F16").Value = 32.71173022246338
$ws.Range("G16").Value = 3.661429423717865
$ws.Range("I16").Value = 25.50018865584435
$ws.Range("K16").Value = 11.41199592553872
$ws.Range("M16").Value = 22.24198505363025
$ws.Range("N16").Value = 20.46902620287952

$ws.Range("B17").Value = 10.91787242428995
$ws.Range("C17").Value = 6.976431217457225
$ws.Range("D17").Value = 6.163276204886613
$ws.Range("F17").Value = 32.63325274295693
$ws.Range("G17").Value = 3.662479482493942
$ws.Range("I17").Value = 25.48801141650016
$ws.Range("K17").Value = 11.34320778619802
$ws.Range("M17").Value = 21.97199665033941
$ws.Range("N17").Value = 20.4868974721746

$ws.Range("B18").Value = 10.86120003117412
$ws.Range("C18").Value = 6.927647488186114
$ws.Range("D18").Value = 6.163575282368998
$ws.Range("F18").Value = 32.58879374316547
$ws.Range("G18").Value = 3.663091508187757
$ws.Range("I18").Value = 25.48147969561
$ws.Range("K18").Value = 11.30372661958873
$ws.Range("M18").Value = 21.81602382265397
$ws.Range("N18").Value = 20.497362399266

$ws.Range("B19").Value = 10.84198992902571
$ws.Range("C19").Value = 6.911082659595989
$ws.Range("D19").Value = 6.163674474729327
$ws.Range("F19").Value = 32.57385823998842
$ws.Range("G19").Value = 3.66330011610463
$ws.Range("I19").Value = 25.47934932453743
$ws.Range("K19").Value = 11.29037477478204
$ws.Range("M19").Value = 21.76310246015861
$ws.Range("N19").Value = 20.50093755767959

$ws.Range("B20").Value = 10.92835032720554
$ws.Range("C20").Value = 6.985436988740969
$ws.Range("D20").Value = 6.163219866775198
$ws.Range("F20").Value = 32.64153672894587
$ws.Range("G20").Value = 3.66236686834285
$ws.Range("I20").Value = 25.48925883225105
$ws.Range("K20").Value = 11.35052204023262
$ws.Range("M20").Value = 22.00080934671277
$ws.Range("N20").Value = 20.48497580948417

$ws.Range("B21").Value = 11.21667910861158
$ws.Range("C21").Value = 7.231694517584033
$ws.Range("D21").Value = 6.161547052442614
$ws.Range("F21").Value = 32.87697314505274
$ws.Range("G21").Value = 3.659327410542731
$ws.Range("I21").Value = 25.5281930069499
$ws.Range("K21").Value = 11.5535101765119
$ws.Range("M21").Value = 22.79091784724586
$ws.Range("N21").Value = 20.43356737878185

$ws.Range("B22").Value = 11.40339835714889
$ws.Range("C22").Value = 7.389695436946169
$ws.Range("D22").Value = 6.160343766864613
$ws.Range("F22").Value = 33.03668084211792
$ws.Range("G22").Value = 3.657412209798853
$ws.Range("I22").Value = 25.55783300246938
$ws.Range("K22").Value = 11.6866167645395
$ws.Range("M22").Value = 23.30004899991037
$ws.Range("N22").Value = 20.40162520971371

$ws.Range("B23").Value = 11.30392235475146
$ws.Range("C23").Value = 7.305654150968609
$ws.Range("D23").Value = 6.160996154396611
$ws.Range("F23").Value = 32.95091840697299
$ws.Range("G23").Value = 3.658427892755105
$ws.Range("I23").Value = 25.54163047869395
$ws.Range("K23").Value = 11.61554942460762
$ws.Range("M23").Value = 23.02903296988406
$ws.Range("N23").Value = 20.41852173634192

$ws.Range("B24").Value = 10.92361377710992
$ws.Range("C24").Value = 6.981366438220062
$ws.Range("D24").Value = 6.163245374490129
$ws.Range("F24").Value = 32.6377894847252
$ws.Range("G24").Value = 3.662417755239828
$ws.Range("I24").Value = 25.48869341422241
$ws.Range("K24").Value = 11.34721505683448
$ws.Range("M24").Value = 21.98778545834662
$ws.Range("N24").Value = 20.48584399984671

$ws.Range("B25").Value = 10.50836305936062
$ws.Range("C25").Value = 6.620869591847392
$ws.Range("D25").Value = 6.165209844591697
$ws.Range("F25").Value = 32.32609972543894
$ws.Range("G25").Value = 3.667029919754631
$ws.Range("I25").Value = 25.44964995341176
$ws.Range("K25").Value = 11.06116926246162
$ws.Range("M25").Value = 20.83938834315423
$ws.Range("N25").Value = 20.56556681083418

